$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the new row 139 from row 138 (keeps header-row style off, applies A/E column styles)
$ws.Range("A138:AC138").Copy() | Out-Null
$ws.Range("A139:AC139").PasteSpecial(-4122) | Out-Null

# Row 4
$ws.Cells.Item(4, 2).Value2 = 6139017
$ws.Cells.Item(4, 7).Value2 = "Harju JK Laagri"
$ws.Cells.Item(4, 8).Value2 = 2
$ws.Cells.Item(4, 9).Value2 = 0
$ws.Cells.Item(4, 10).Value2 = "H"
$ws.Cells.Item(4, 11).Value2 = 1.666
$ws.Cells.Item(4, 12).Value2 = 3.6
$ws.Cells.Item(4, 13).Value2 = 4.2
$ws.Cells.Item(4, 14).Value2 = 1.727
$ws.Cells.Item(4, 15).Value2 = 3.5
$ws.Cells.Item(4, 16).Value2 = 4
$ws.Cells.Item(4, 17).Value2 = -0.75
$ws.Cells.Item(4, 18).Value2 = 2
$ws.Cells.Item(4, 19).Value2 = 1.8
$ws.Cells.Item(4, 20).Value2 = 2.5
$ws.Cells.Item(4, 21).Value2 = 1.9
$ws.Cells.Item(4, 22).Value2 = 1.9
$ws.Cells.Item(4, 23).Value2 = 0.7270000000000001
$ws.Cells.Item(4, 25).Value2 = -1
$ws.Cells.Item(4, 26).Value2 = 1
$ws.Cells.Item(4, 27).Value2 = -1
$ws.Cells.Item(4, 29).Value2 = 0.8999999999999999

# Row 5
$ws.Cells.Item(5, 2).Value2 = 6139018
$ws.Cells.Item(5, 7).Value2 = "JK Trans Narva"
$ws.Cells.Item(5, 8).Value2 = 0
$ws.Cells.Item(5, 9).Value2 = 1
$ws.Cells.Item(5, 10).Value2 = "A"
$ws.Cells.Item(5, 11).Value2 = 2.4
$ws.Cells.Item(5, 12).Value2 = 3.4
$ws.Cells.Item(5, 13).Value2 = 2.5
$ws.Cells.Item(5, 14).Value2 = 2.875
$ws.Cells.Item(5, 15).Value2 = 3.1
$ws.Cells.Item(5, 16).Value2 = 2.3
$ws.Cells.Item(5, 17).Value2 = 0.25
$ws.Cells.Item(5, 18).Value2 = 1.75
$ws.Cells.Item(5, 19).Value2 = 2.05
$ws.Cells.Item(5, 20).Value2 = 2.25
$ws.Cells.Item(5, 21).Value2 = 1.925
$ws.Cells.Item(5, 22).Value2 = 1.875
$ws.Cells.Item(5, 23).Value2 = -1
$ws.Cells.Item(5, 25).Value2 = 1.3
$ws.Cells.Item(5, 26).Value2 = -1
$ws.Cells.Item(5, 27).Value2 = 1.05
$ws.Cells.Item(5, 29).Value2 = 0.875

# Row 10
$ws.Cells.Item(10, 7).Value2 = "JK Tammeka Tartu"

# Row 11
$ws.Cells.Item(11, 7).Value2 = "JK Tallinna Kalev"

# Row 13
$ws.Cells.Item(13, 6).Value2 = "JK Tallinna Kalev"

# Row 15
$ws.Cells.Item(15, 6).Value2 = "JK Tammeka Tartu"

# Row 20
$ws.Cells.Item(20, 7).Value2 = "JK Tammeka Tartu"

# Row 21
$ws.Cells.Item(21, 7).Value2 = "JK Tallinna Kalev"

# Row 24
$ws.Cells.Item(24, 6).Value2 = "JK Tammeka Tartu"

# Row 25
$ws.Cells.Item(25, 7).Value2 = "JK Tallinna Kalev"

# Row 27
$ws.Cells.Item(27, 6).Value2 = "JK Tammeka Tartu"

# Row 30
$ws.Cells.Item(30, 6).Value2 = "JK Tallinna Kalev"

# Row 33
$ws.Cells.Item(33, 7).Value2 = "JK Tallinna Kalev"

# Row 36
$ws.Cells.Item(36, 6).Value2 = "JK Tammeka Tartu"
$ws.Cells.Item(36, 7).Value2 = "JK Tallinna Kalev"

# Row 39
$ws.Cells.Item(39, 7).Value2 = "JK Tammeka Tartu"

# Row 41
$ws.Cells.Item(41, 6).Value2 = "JK Tallinna Kalev"

# Row 44
$ws.Cells.Item(44, 7).Value2 = "JK Tallinna Kalev"

# Row 47
$ws.Cells.Item(47, 6).Value2 = "JK Tammeka Tartu"

# Row 49
$ws.Cells.Item(49, 6).Value2 = "JK Tallinna Kalev"

# Row 50
$ws.Cells.Item(50, 7).Value2 = "JK Tammeka Tartu"

# Row 53
$ws.Cells.Item(53, 7).Value2 = "JK Tammeka Tartu"

# Row 55
$ws.Cells.Item(55, 7).Value2 = "JK Tallinna Kalev"

# Row 58
$ws.Cells.Item(58, 6).Value2 = "JK Tallinna Kalev"

# Row 59
$ws.Cells.Item(59, 6).Value2 = "JK Tammeka Tartu"

# Row 63
$ws.Cells.Item(63, 7).Value2 = "JK Tammeka Tartu"

# Row 66
$ws.Cells.Item(66, 7).Value2 = "JK Tallinna Kalev"

# Row 70
$ws.Cells.Item(70, 6).Value2 = "JK Tammeka Tartu"

# Row 71
$ws.Cells.Item(71, 2).Value2 = 6139071
$ws.Cells.Item(71, 6).Value2 = "Parnu JK Vaprus"
$ws.Cells.Item(71, 7).Value2 = "JK Trans Narva"
$ws.Cells.Item(71, 8).Value2 = 3
$ws.Cells.Item(71, 10).Value2 = "H"
$ws.Cells.Item(71, 11).Value2 = 2.4
$ws.Cells.Item(71, 12).Value2 = 3.2
$ws.Cells.Item(71, 13).Value2 = 2.6
$ws.Cells.Item(71, 14).Value2 = 3
$ws.Cells.Item(71, 15).Value2 = 3.25
$ws.Cells.Item(71, 16).Value2 = 2.2
$ws.Cells.Item(71, 17).Value2 = 0.25
$ws.Cells.Item(71, 18).Value2 = 1.825
$ws.Cells.Item(71, 19).Value2 = 1.975
$ws.Cells.Item(71, 20).Value2 = 2.5
$ws.Cells.Item(71, 21).Value2 = 1.875
$ws.Cells.Item(71, 22).Value2 = 1.925
$ws.Cells.Item(71, 23).Value2 = 2
$ws.Cells.Item(71, 25).Value2 = -1
$ws.Cells.Item(71, 26).Value2 = 0.825
$ws.Cells.Item(71, 28).Value2 = 0.875
$ws.Cells.Item(71, 29).Value2 = -1

# Row 72
$ws.Cells.Item(72, 2).Value2 = 6139072
$ws.Cells.Item(72, 6).Value2 = "JK Tallinna Kalev"
$ws.Cells.Item(72, 7).Value2 = "FC Flora Tallinn"
$ws.Cells.Item(72, 8).Value2 = 1
$ws.Cells.Item(72, 10).Value2 = "A"
$ws.Cells.Item(72, 11).Value2 = 9
$ws.Cells.Item(72, 12).Value2 = 7
$ws.Cells.Item(72, 13).Value2 = 1.166
$ws.Cells.Item(72, 14).Value2 = 7
$ws.Cells.Item(72, 15).Value2 = 6
$ws.Cells.Item(72, 16).Value2 = 1.25
$ws.Cells.Item(72, 17).Value2 = 1.75
$ws.Cells.Item(72, 18).Value2 = 1.9
$ws.Cells.Item(72, 19).Value2 = 1.9
$ws.Cells.Item(72, 20).Value2 = 3
$ws.Cells.Item(72, 21).Value2 = 1.95
$ws.Cells.Item(72, 22).Value2 = 1.85
$ws.Cells.Item(72, 23).Value2 = -1
$ws.Cells.Item(72, 25).Value2 = 0.25
$ws.Cells.Item(72, 26).Value2 = 0.8999999999999999
$ws.Cells.Item(72, 28).Value2 = 0
$ws.Cells.Item(72, 29).Value2 = -0

# Row 74
$ws.Cells.Item(74, 6).Value2 = "JK Tallinna Kalev"

# Row 75
$ws.Cells.Item(75, 7).Value2 = "JK Tammeka Tartu"

# Row 79
$ws.Cells.Item(79, 6).Value2 = "JK Tammeka Tartu"

# Row 80
$ws.Cells.Item(80, 7).Value2 = "JK Tallinna Kalev"

# Row 83
$ws.Cells.Item(83, 6).Value2 = "JK Tammeka Tartu"

# Row 85
$ws.Cells.Item(85, 7).Value2 = "JK Tallinna Kalev"

# Row 88
$ws.Cells.Item(88, 2).Value2 = 6376947
$ws.Cells.Item(88, 6).Value2 = "JK Tallinna Kalev"
$ws.Cells.Item(88, 7).Value2 = "JK Tammeka Tartu"
$ws.Cells.Item(88, 8).Value2 = 2
$ws.Cells.Item(88, 9).Value2 = 7
$ws.Cells.Item(88, 10).Value2 = "A"
$ws.Cells.Item(88, 11).Value2 = 3.6
$ws.Cells.Item(88, 12).Value2 = 3.4
$ws.Cells.Item(88, 13).Value2 = 1.909
$ws.Cells.Item(88, 14).Value2 = 2.4
$ws.Cells.Item(88, 15).Value2 = 3.6
$ws.Cells.Item(88, 16).Value2 = 2.45
$ws.Cells.Item(88, 17).Value2 = 0
$ws.Cells.Item(88, 20).Value2 = 2.75
$ws.Cells.Item(88, 21).Value2 = 1.975
$ws.Cells.Item(88, 22).Value2 = 1.825
$ws.Cells.Item(88, 24).Value2 = -1
$ws.Cells.Item(88, 25).Value2 = 1.45
$ws.Cells.Item(88, 28).Value2 = 0.9750000000000001
$ws.Cells.Item(88, 29).Value2 = -1

# Row 89
$ws.Cells.Item(89, 2).Value2 = 6376945
$ws.Cells.Item(89, 6).Value2 = "Parnu JK Vaprus"
$ws.Cells.Item(89, 7).Value2 = "Harju JK Laagri"
$ws.Cells.Item(89, 8).Value2 = 0
$ws.Cells.Item(89, 9).Value2 = 0
$ws.Cells.Item(89, 10).Value2 = "D"
$ws.Cells.Item(89, 11).Value2 = 1.615
$ws.Cells.Item(89, 12).Value2 = 4
$ws.Cells.Item(89, 13).Value2 = 4.5
$ws.Cells.Item(89, 14).Value2 = 1.85
$ws.Cells.Item(89, 15).Value2 = 3.8
$ws.Cells.Item(89, 16).Value2 = 3.5
$ws.Cells.Item(89, 17).Value2 = -0.5
$ws.Cells.Item(89, 20).Value2 = 2.5
$ws.Cells.Item(89, 21).Value2 = 1.75
$ws.Cells.Item(89, 22).Value2 = 1.95
$ws.Cells.Item(89, 24).Value2 = 2.8
$ws.Cells.Item(89, 25).Value2 = -1
$ws.Cells.Item(89, 28).Value2 = -1
$ws.Cells.Item(89, 29).Value2 = 0.95

# Row 93
$ws.Cells.Item(93, 7).Value2 = "JK Tammeka Tartu"

# Row 95
$ws.Cells.Item(95, 6).Value2 = "JK Tallinna Kalev"

# Row 100
$ws.Cells.Item(100, 6).Value2 = "JK Tammeka Tartu"

# Row 102
$ws.Cells.Item(102, 7).Value2 = "JK Tallinna Kalev"

# Row 103
$ws.Cells.Item(103, 6).Value2 = "JK Tallinna Kalev"

# Row 104
$ws.Cells.Item(104, 2).Value2 = 6533597
$ws.Cells.Item(104, 6).Value2 = "FC Kuressaare"
$ws.Cells.Item(104, 7).Value2 = "Parnu JK Vaprus"
$ws.Cells.Item(104, 8).Value2 = 1
$ws.Cells.Item(104, 11).Value2 = 2.5
$ws.Cells.Item(104, 12).Value2 = 3.4
$ws.Cells.Item(104, 13).Value2 = 2.5
$ws.Cells.Item(104, 14).Value2 = 2.15
$ws.Cells.Item(104, 15).Value2 = 3.6
$ws.Cells.Item(104, 16).Value2 = 2.875
$ws.Cells.Item(104, 17).Value2 = -0.25
$ws.Cells.Item(104, 18).Value2 = 1.95
$ws.Cells.Item(104, 19).Value2 = 1.85
$ws.Cells.Item(104, 21).Value2 = 1.95
$ws.Cells.Item(104, 22).Value2 = 1.85
$ws.Cells.Item(104, 23).Value2 = 1.15
$ws.Cells.Item(104, 26).Value2 = 0.95
$ws.Cells.Item(104, 28).Value2 = -1
$ws.Cells.Item(104, 29).Value2 = 0.8500000000000001

# Row 105
$ws.Cells.Item(105, 2).Value2 = 6537957
$ws.Cells.Item(105, 6).Value2 = "FC Flora Tallinn"
$ws.Cells.Item(105, 7).Value2 = "JK Nomme Kalju"
$ws.Cells.Item(105, 8).Value2 = 0
$ws.Cells.Item(105, 9).Value2 = 0
$ws.Cells.Item(105, 11).Value2 = 1.4
$ws.Cells.Item(105, 12).Value2 = 4
$ws.Cells.Item(105, 13).Value2 = 7.5
$ws.Cells.Item(105, 14).Value2 = 1.5
$ws.Cells.Item(105, 15).Value2 = 4.2
$ws.Cells.Item(105, 16).Value2 = 5
$ws.Cells.Item(105, 17).Value2 = -1
$ws.Cells.Item(105, 21).Value2 = 1.85
$ws.Cells.Item(105, 22).Value2 = 1.95
$ws.Cells.Item(105, 24).Value2 = 3.2
$ws.Cells.Item(105, 26).Value2 = -1
$ws.Cells.Item(105, 27).Value2 = 0.95
$ws.Cells.Item(105, 28).Value2 = -1
$ws.Cells.Item(105, 29).Value2 = 0.95

# Row 106
$ws.Cells.Item(106, 2).Value2 = 6535416
$ws.Cells.Item(106, 6).Value2 = "Paide Linnameeskond"
$ws.Cells.Item(106, 7).Value2 = "FC Levadia Tallinn"
$ws.Cells.Item(106, 8).Value2 = 2
$ws.Cells.Item(106, 9).Value2 = 2
$ws.Cells.Item(106, 10).Value2 = "D"
$ws.Cells.Item(106, 11).Value2 = 3
$ws.Cells.Item(106, 12).Value2 = 3.8
$ws.Cells.Item(106, 13).Value2 = 2
$ws.Cells.Item(106, 14).Value2 = 3
$ws.Cells.Item(106, 15).Value2 = 4
$ws.Cells.Item(106, 16).Value2 = 1.909
$ws.Cells.Item(106, 17).Value2 = 0.5
$ws.Cells.Item(106, 18).Value2 = 1.85
$ws.Cells.Item(106, 19).Value2 = 1.95
$ws.Cells.Item(106, 23).Value2 = -1
$ws.Cells.Item(106, 24).Value2 = 3
$ws.Cells.Item(106, 26).Value2 = 0.8500000000000001
$ws.Cells.Item(106, 28).Value2 = 0.95
$ws.Cells.Item(106, 29).Value2 = -1

# Row 107
$ws.Cells.Item(107, 2).Value2 = 6537869
$ws.Cells.Item(107, 6).Value2 = "JK Tammeka Tartu"
$ws.Cells.Item(107, 7).Value2 = "JK Trans Narva"
$ws.Cells.Item(107, 8).Value2 = 5
$ws.Cells.Item(107, 10).Value2 = "H"
$ws.Cells.Item(107, 11).Value2 = 1.6
$ws.Cells.Item(107, 13).Value2 = 4.5
$ws.Cells.Item(107, 14).Value2 = 1.65
$ws.Cells.Item(107, 15).Value2 = 4
$ws.Cells.Item(107, 16).Value2 = 4.333
$ws.Cells.Item(107, 17).Value2 = -0.75
$ws.Cells.Item(107, 18).Value2 = 1.8
$ws.Cells.Item(107, 19).Value2 = 2
$ws.Cells.Item(107, 21).Value2 = 1.9
$ws.Cells.Item(107, 22).Value2 = 1.9
$ws.Cells.Item(107, 23).Value2 = 0.6499999999999999
$ws.Cells.Item(107, 24).Value2 = -1
$ws.Cells.Item(107, 26).Value2 = 0.8
$ws.Cells.Item(107, 27).Value2 = -1
$ws.Cells.Item(107, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(107, 29).Value2 = -1

# Row 108
$ws.Cells.Item(108, 6).Value2 = "JK Tallinna Kalev"

# Row 112
$ws.Cells.Item(112, 6).Value2 = "JK Tammeka Tartu"
$ws.Cells.Item(112, 7).Value2 = "JK Tallinna Kalev"

# Row 114
$ws.Cells.Item(114, 6).Value2 = "JK Tallinna Kalev"

# Row 115
$ws.Cells.Item(115, 2).Value2 = 7919323
$ws.Cells.Item(115, 6).Value2 = "JK Nomme Kalju"
$ws.Cells.Item(115, 7).Value2 = "JK Trans Narva"
$ws.Cells.Item(115, 8).Value2 = 3
$ws.Cells.Item(115, 9).Value2 = 0
$ws.Cells.Item(115, 10).Value2 = "H"
$ws.Cells.Item(115, 11).Value2 = 1.285
$ws.Cells.Item(115, 12).Value2 = 5.5
$ws.Cells.Item(115, 13).Value2 = 6.5
$ws.Cells.Item(115, 14).Value2 = 1.571
$ws.Cells.Item(115, 15).Value2 = 4.75
$ws.Cells.Item(115, 16).Value2 = 4.2
$ws.Cells.Item(115, 17).Value2 = -1
$ws.Cells.Item(115, 18).Value2 = 1.925
$ws.Cells.Item(115, 19).Value2 = 1.875
$ws.Cells.Item(115, 20).Value2 = 2.75
$ws.Cells.Item(115, 21).Value2 = 1.875
$ws.Cells.Item(115, 22).Value2 = 1.925
$ws.Cells.Item(115, 23).Value2 = 0.571
$ws.Cells.Item(115, 25).Value2 = -1
$ws.Cells.Item(115, 26).Value2 = 0.925
$ws.Cells.Item(115, 27).Value2 = -1
$ws.Cells.Item(115, 28).Value2 = 0.4375
$ws.Cells.Item(115, 29).Value2 = -0.5

# Row 116
$ws.Cells.Item(116, 2).Value2 = 7919322
$ws.Cells.Item(116, 6).Value2 = "FC Kuressaare"
$ws.Cells.Item(116, 7).Value2 = "FC Levadia Tallinn"
$ws.Cells.Item(116, 8).Value2 = 0
$ws.Cells.Item(116, 9).Value2 = 6
$ws.Cells.Item(116, 10).Value2 = "A"
$ws.Cells.Item(116, 11).Value2 = 11
$ws.Cells.Item(116, 12).Value2 = 6
$ws.Cells.Item(116, 13).Value2 = 1.166
$ws.Cells.Item(116, 14).Value2 = 15
$ws.Cells.Item(116, 15).Value2 = 8.5
$ws.Cells.Item(116, 16).Value2 = 1.125
$ws.Cells.Item(116, 17).Value2 = 2.5
$ws.Cells.Item(116, 18).Value2 = 1.825
$ws.Cells.Item(116, 19).Value2 = 1.975
$ws.Cells.Item(116, 20).Value2 = 3.25
$ws.Cells.Item(116, 21).Value2 = 1.9
$ws.Cells.Item(116, 22).Value2 = 1.9
$ws.Cells.Item(116, 23).Value2 = -1
$ws.Cells.Item(116, 25).Value2 = 0.125
$ws.Cells.Item(116, 26).Value2 = -1
$ws.Cells.Item(116, 27).Value2 = 0.9750000000000001
$ws.Cells.Item(116, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(116, 29).Value2 = -1

# Row 118
$ws.Cells.Item(118, 7).Value2 = "JK Tammeka Tartu"

# Row 119
$ws.Cells.Item(119, 6).Value2 = "JK Tammeka Tartu"

# Row 121
$ws.Cells.Item(121, 7).Value2 = "JK Tallinna Kalev"

# Row 124
$ws.Cells.Item(124, 7).Value2 = "JK Tammeka Tartu"

# Row 127
$ws.Cells.Item(127, 6).Value2 = "JK Tallinna Kalev"

# Row 129
$ws.Cells.Item(129, 7).Value2 = "JK Tallinna Kalev"

# Row 133
$ws.Cells.Item(133, 6).Value2 = "JK Tammeka Tartu"

# Row 134
$ws.Cells.Item(134, 6).Value2 = "JK Tallinna Kalev"

# Row 135
$ws.Cells.Item(135, 7).Value2 = "JK Tammeka Tartu"

# Row 139
$ws.Cells.Item(139, 1).Value2 = 137
$ws.Cells.Item(139, 2).Value2 = 7719670
$ws.Cells.Item(139, 3).Value2 = "Estonia Meistriliiga"
$ws.Cells.Item(139, 4).Value2 = "Estonia Meistriliiga"
$ws.Cells.Item(139, 5).Value2 = 45398.54166666666
$ws.Cells.Item(139, 6).Value2 = "JK Tallinna Kalev"
$ws.Cells.Item(139, 7).Value2 = "FC Levadia Tallinn"
$ws.Cells.Item(139, 8).Value2 = 1
$ws.Cells.Item(139, 9).Value2 = 2
$ws.Cells.Item(139, 10).Value2 = "A"
$ws.Cells.Item(139, 11).Value2 = 7
$ws.Cells.Item(139, 12).Value2 = 5
$ws.Cells.Item(139, 13).Value2 = 1.3
$ws.Cells.Item(139, 14).Value2 = 29
$ws.Cells.Item(139, 15).Value2 = 12
$ws.Cells.Item(139, 16).Value2 = 1.055
$ws.Cells.Item(139, 17).Value2 = 3
$ws.Cells.Item(139, 18).Value2 = 1.95
$ws.Cells.Item(139, 19).Value2 = 1.85
$ws.Cells.Item(139, 20).Value2 = 3.5
$ws.Cells.Item(139, 21).Value2 = 1.775
$ws.Cells.Item(139, 22).Value2 = 1.925
$ws.Cells.Item(139, 23).Value2 = -1
$ws.Cells.Item(139, 24).Value2 = -1
$ws.Cells.Item(139, 25).Value2 = 0.05499999999999994
$ws.Cells.Item(139, 26).Value2 = 0.95
$ws.Cells.Item(139, 27).Value2 = -1
$ws.Cells.Item(139, 28).Value2 = -1
$ws.Cells.Item(139, 29).Value2 = 0.925
